# Apply updated cryptos list values (generated from target diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-detected as a number by Excel;
# force them to stay Text (matching the source inlineStr cells) via NumberFormat.
$ws.Range("D2").Value = "26.276.58"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "1.618.54"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.92"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.484"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.248"
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.73"
$ws.Range("E10").Value = "  +5.12%  "
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").Value = "1.843.28"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").Value = "1.623.32"
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.00"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").Value = "26.291.39"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.22"
$ws.Range("E17").Value = "  +3.91%  "
$ws.Range("D18").Value = "0.0₃0725"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "201.30"
$ws.Range("E20").Value = "  +1.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.28"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.04"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.89"
$ws.Range("E24").Value = "  +2.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.43"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.15"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.54"
$ws.Range("E29").Value = "  +1.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0518"
$ws.Range("E30").Value = "  +9.79%  "
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.49"
$ws.Range("E34").Value = "  +2.11%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.40"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("D36").Value = "1.176.41"
$ws.Range("E36").Value = "  +5.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0163"
$ws.Range("E37").Value = "  +2.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.804"
$ws.Range("E38").Value = "  +3.37%  "
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.32"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.789"
$ws.Range("E42").Value = "  +1.48%  "
$ws.Range("E43").Value = "  +5.01%  "
$ws.Range("D44").Value = "1.754.30"
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("E46").Value = "  +14.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.53"
$ws.Range("E47").Value = "  +4.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.61"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.409"
$ws.Range("E50").Value = "  +0.57%  "
